{"js": "// Replace the date line and the 25 division-problem answers in the table\n// with their new values, per the diff. Every old string is unique within\n// the document, so a direct search + replace for each pair is safe.\nconst replacements = [\n  [\"2025-09-29 Monday\", \"2025-09-30 Tuesday\"],\n  [\"717\\u00F72=358, 1\", \"463\\u00F77=66, 1\"],\n  [\"431\\u00F75=86, 1\", \"708\\u00F75=141, 3\"],\n  [\"566\\u00F77=80, 6\", \"414\\u00F77=59, 1\"],\n  [\"914\\u00F77=130, 4\", \"173\\u00F79=19, 2\"],\n  [\"264\\u00F77=37, 5\", \"324\\u00F78=40, 4\"],\n  [\"874\\u00F77=124, 6\", \"429\\u00F76=71, 3\"],\n  [\"212\\u00F76=35, 2\", \"526\\u00F75=105, 1\"],\n  [\"387\\u00F79=43, 0\", \"455\\u00F75=91, 0\"],\n  [\"648\\u00F76=108, 0\", \"437\\u00F77=62, 3\"],\n  [\"402\\u00F79=44, 6\", \"489\\u00F72=244, 1\"],\n  [\"651\\u00F72=325, 1\", \"883\\u00F72=441, 1\"],\n  [\"995\\u00F74=248, 3\", \"340\\u00F75=68, 0\"],\n  [\"781\\u00F76=130, 1\", \"580\\u00F73=193, 1\"],\n  [\"698\\u00F75=139, 3\", \"310\\u00F72=155, 0\"],\n  [\"157\\u00F74=39, 1\", \"647\\u00F78=80, 7\"],\n  [\"695\\u00F79=77, 2\", \"519\\u00F75=103, 4\"],\n  [\"870\\u00F77=124, 2\", \"123\\u00F74=30, 3\"],\n  [\"461\\u00F72=230, 1\", \"499\\u00F79=55, 4\"],\n  [\"857\\u00F72=428, 1\", \"548\\u00F72=274, 0\"],\n  [\"639\\u00F76=106, 3\", \"833\\u00F72=416, 1\"],\n  [\"694\\u00F77=99, 1\", \"787\\u00F79=87, 4\"],\n  [\"279\\u00F75=55, 4\", \"894\\u00F72=447, 0\"],\n  [\"540\\u00F75=108, 0\", \"173\\u00F74=43, 1\"],\n  [\"282\\u00F78=35, 2\", \"519\\u00F74=129, 3\"],\n  [\"830\\u00F77=118, 4\", \"761\\u00F75=152, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the 25 division-problem answers in the table\n# with their new values, per the diff. Every old string is unique within\n# the document, so a direct Find/Replace for each pair is safe.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-09-29 Monday\", \"2025-09-30 Tuesday\"),\n    @(\"717\u00f72=358, 1\", \"463\u00f77=66, 1\"),\n    @(\"431\u00f75=86, 1\", \"708\u00f75=141, 3\"),\n    @(\"566\u00f77=80, 6\", \"414\u00f77=59, 1\"),\n    @(\"914\u00f77=130, 4\", \"173\u00f79=19, 2\"),\n    @(\"264\u00f77=37, 5\", \"324\u00f78=40, 4\"),\n    @(\"874\u00f77=124, 6\", \"429\u00f76=71, 3\"),\n    @(\"212\u00f76=35, 2\", \"526\u00f75=105, 1\"),\n    @(\"387\u00f79=43, 0\", \"455\u00f75=91, 0\"),\n    @(\"648\u00f76=108, 0\", \"437\u00f77=62, 3\"),\n    @(\"402\u00f79=44, 6\", \"489\u00f72=244, 1\"),\n    @(\"651\u00f72=325, 1\", \"883\u00f72=441, 1\"),\n    @(\"995\u00f74=248, 3\", \"340\u00f75=68, 0\"),\n    @(\"781\u00f76=130, 1\", \"580\u00f73=193, 1\"),\n    @(\"698\u00f75=139, 3\", \"310\u00f72=155, 0\"),\n    @(\"157\u00f74=39, 1\", \"647\u00f78=80, 7\"),\n    @(\"695\u00f79=77, 2\", \"519\u00f75=103, 4\"),\n    @(\"870\u00f77=124, 2\", \"123\u00f74=30, 3\"),\n    @(\"461\u00f72=230, 1\", \"499\u00f79=55, 4\"),\n    @(\"857\u00f72=428, 1\", \"548\u00f72=274, 0\"),\n    @(\"639\u00f76=106, 3\", \"833\u00f72=416, 1\"),\n    @(\"694\u00f77=99, 1\", \"787\u00f79=87, 4\"),\n    @(\"279\u00f75=55, 4\", \"894\u00f72=447, 0\"),\n    @(\"540\u00f75=108, 0\", \"173\u00f74=43, 1\"),\n    @(\"282\u00f78=35, 2\", \"519\u00f74=129, 3\"),\n    @(\"830\u00f77=118, 4\", \"761\u00f75=152, 1\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
